$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.760.95"
$ws.Range("E2").Value = "'  +2.33%  "
$ws.Range("D3").Value = "'3.033.87"
$ws.Range("E3").Value = "'  +1.94%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'512.02"
$ws.Range("E5").Value = "'  +1.91%  "
$ws.Range("D6").Value = "'140.11"
$ws.Range("E6").Value = "'  +3.91%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'0.443"
$ws.Range("E8").Value = "'  +3.35%  "
$ws.Range("D9").Value = "'7.51"
$ws.Range("E9").Value = "'  +0.81%  "
$ws.Range("E10").Value = "'  +3.68%  "
$ws.Range("E11").Value = "'  +5.33%  "
$ws.Range("D12").Value = "'3.549.64"
$ws.Range("E12").Value = "'  +1.92%  "
$ws.Range("E13").Value = "'  +2.04%  "
$ws.Range("E14").Value = "'  +5.56%  "
$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "'  +10.51%  "
$ws.Range("D16").Value = "'57.745.73"
$ws.Range("E16").Value = "'  +2.38%  "
$ws.Range("D17").Value = "'6.25"
$ws.Range("E17").Value = "'  +9.07%  "
$ws.Range("D18").Value = "'3.034.57"
$ws.Range("E18").Value = "'  +2.13%  "
$ws.Range("E19").Value = "'  +4.97%  "
$ws.Range("E20").Value = "'  +3.79%  "
$ws.Range("D21").Value = "'334.01"
$ws.Range("E21").Value = "'  +3.80%  "
$ws.Range("B22").Value = "'Dai"
$ws.Range("C22").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("B23").Value = "'LEO"
$ws.Range("C23").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'5.76"
$ws.Range("E23").Value = "'  +1.07%  "
$ws.Range("E24").Value = "'  +6.87%  "
$ws.Range("D25").Value = "'64.67"
$ws.Range("E25").Value = "'  +4.80%  "
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "'  +5.29%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  +0.07%  "
$ws.Range("D28").Value = "'0.0₃0932"
$ws.Range("E28").Value = "'  +5.18%  "
$ws.Range("D29").Value = "'6.85"
$ws.Range("E29").Value = "'  +6.97%  "
$ws.Range("D30").Value = "'7.50"
$ws.Range("E30").Value = "'  +11.17%  "
$ws.Range("E31").Value = "'  +3.52%  "
$ws.Range("E32").Value = "'  +3.83%  "
$ws.Range("E33").Value = "'  +2.64%  "
$ws.Range("D34").Value = "'155.90"
$ws.Range("E34").Value = "'  -1.58%  "
$ws.Range("E35").Value = "'  +6.53%  "
$ws.Range("E36").Value = "'  +6.39%  "
$ws.Range("E37").Value = "'  +2.31%  "
$ws.Range("D38").Value = "'24.88"
$ws.Range("E38").Value = "'  +8.55%  "
$ws.Range("D39").Value = "'0.0688"
$ws.Range("E39").Value = "'  +2.49%  "
$ws.Range("D40").Value = "'3.066.53"
$ws.Range("E40").Value = "'  +1.91%  "
$ws.Range("D41").Value = "'37.44"
$ws.Range("E41").Value = "'  +3.41%  "
$ws.Range("E42").Value = "'  +9.60%  "
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("D44").Value = "'2.313.54"
$ws.Range("E44").Value = "'  +3.22%  "
$ws.Range("D45").Value = "'0.657"
$ws.Range("E45").Value = "'  +2.87%  "
$ws.Range("E46").Value = "'  +2.68%  "
$ws.Range("D47").Value = "'0.996"
$ws.Range("E47").Value = "'  +1.65%  "
$ws.Range("D48").Value = "'6.04"
$ws.Range("E48").Value = "'  +5.48%  "
$ws.Range("E49").Value = "'  +2.73%  "
$ws.Range("D50").Value = "'19.85"
$ws.Range("E50").Value = "'  +4.91%  "
$ws.Range("D51").Value = "'1.85"
$ws.Range("E51").Value = "'  -3.75%  "
